$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 284.875
$ws.Range("I4").Value = 284.875
$ws.Range("K4").Value = 284.875
$ws.Range("M4").Value = -170.875

$ws.Range("H6").Value = 243.4
$ws.Range("I6").Value = 243.4
$ws.Range("K6").Value = 730.2
$ws.Range("M6").Value = -618.2

$ws.Range("H28").Value = 11278.944
$ws.Range("I28").Value = 3720.5833
$ws.Range("J28").Value = 26395.666
$ws.Range("K28").Value = 3720.5833
$ws.Range("L28").Value = 26395.666
$ws.Range("M28").Value = -3235.5833
$ws.Range("N28").Value = -27365.666

$ws.Range("H76").Value = 3888
$ws.Range("I76").Value = 3888
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 3888
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -3573
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 3888
$ws.Range("I79").Value = 3888
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 3888
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -2796
$ws.Range("N79").ClearContents()

$ws.Range("H112").Value = 4151.143
$ws.Range("I112").Value = 1746
$ws.Range("J112").Value = 5113.2
$ws.Range("K112").Value = 5238
$ws.Range("L112").Value = 15339.6
$ws.Range("M112").Value = -4130
$ws.Range("N112").Value = -17555.6

$ws.Range("H132").Value = 185363.81
$ws.Range("I132").Value = 4144.2856
$ws.Range("J132").Value = 502498
$ws.Range("K132").Value = 12432.8568
$ws.Range("L132").Value = 1507494
$ws.Range("M132").Value = -9902.856800000001
$ws.Range("N132").Value = -1512554

$ws.Range("H137").Value = 2119.7742
$ws.Range("I137").Value = 1163.9333
$ws.Range("J137").Value = 3015.875
$ws.Range("K137").Value = 3491.7999
$ws.Range("L137").Value = 9047.625
$ws.Range("M137").Value = -941.7999
$ws.Range("N137").Value = -14147.625

$ws.Range("H138").Value = 9548.143
$ws.Range("I138").Value = 6995
$ws.Range("J138").Value = 9973.666999999999
$ws.Range("K138").Value = 20985
$ws.Range("L138").Value = 29921.001
$ws.Range("M138").Value = -15845
$ws.Range("N138").Value = -40201.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4804.0557
$ws.Range("I32").Value = 2698.28
$ws.Range("K32").Value = 2698.28
$ws.Range("M32").Value = -2411.28

$ws.Range("H45").Value = 1965.5
$ws.Range("I45").Value = 1380.1578
$ws.Range("J45").Value = 3554.2856
$ws.Range("K45").Value = 1380.1578
$ws.Range("L45").Value = 3554.2856
$ws.Range("M45").Value = -1003.1578
$ws.Range("N45").Value = -4308.2856

$ws.Range("H61").Value = 8899.4
$ws.Range("I61").Value = 26000
$ws.Range("K61").Value = 26000
$ws.Range("M61").Value = -25788

$ws.Range("H63").Value = 1952.7222
$ws.Range("I63").Value = 1175.7858
$ws.Range("J63").Value = 4672
$ws.Range("K63").Value = 1175.7858
$ws.Range("L63").Value = 4672
$ws.Range("M63").Value = -489.7858000000001
$ws.Range("N63").Value = -6044

$ws.Range("H66").Value = 1952.7222
$ws.Range("I66").Value = 1175.7858
$ws.Range("J66").Value = 4672
$ws.Range("K66").Value = 5878.929
$ws.Range("L66").Value = 23360
$ws.Range("M66").Value = -2446.929
$ws.Range("N66").Value = -30224

$ws.Range("H122").Value = 1242.4546
$ws.Range("I122").Value = 1296.9
$ws.Range("J122").Value = 698
$ws.Range("K122").Value = 3890.7
$ws.Range("L122").Value = 2094
$ws.Range("M122").Value = -1440.7
$ws.Range("N122").Value = -6994

$ws.Range("H136").Value = 8899.4
$ws.Range("I136").Value = 26000
$ws.Range("K136").Value = 78000
$ws.Range("M136").Value = -75450

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1997.8572
$ws.Range("I99").Value = 1985
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 1985
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -487
$ws.Range("N99").Value = -4996

$ws.Range("H105").Value = 3195.3333
$ws.Range("I105").Value = 3156
$ws.Range("K105").Value = 3156
$ws.Range("M105").Value = -1409

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1385.8572
$ws.Range("J22").Value = 1450.1666
$ws.Range("L22").Value = 1450.1666
$ws.Range("N22").Value = -2150.1666

$ws.Range("H31").Value = 4211.8076
$ws.Range("I31").Value = 1574.6
$ws.Range("J31").Value = 4839.7144
$ws.Range("K31").Value = 1574.6
$ws.Range("L31").Value = 4839.7144
$ws.Range("M31").Value = -1279.6
$ws.Range("N31").Value = -5429.7144

$ws.Range("H34").Value = 4211.8076
$ws.Range("I34").Value = 1574.6
$ws.Range("J34").Value = 4839.7144
$ws.Range("K34").Value = 1574.6
$ws.Range("L34").Value = 4839.7144
$ws.Range("M34").Value = -1372.6
$ws.Range("N34").Value = -5243.7144

$ws.Range("H58").Value = 6091.2104
$ws.Range("I58").Value = 6229.4165
$ws.Range("J58").Value = 5854.2856
$ws.Range("K58").Value = 6229.4165
$ws.Range("L58").Value = 5854.2856
$ws.Range("M58").Value = -6026.4165
$ws.Range("N58").Value = -6260.2856

$ws.Range("H99").Value = 19303.273
$ws.Range("I99").Value = 4562.4
$ws.Range("J99").Value = 31587.334
$ws.Range("K99").Value = 4562.4
$ws.Range("L99").Value = 31587.334
$ws.Range("M99").Value = -3064.4
$ws.Range("N99").Value = -34583.334

$ws.Range("H122").Value = 4192.7144
$ws.Range("I122").Value = 4200
$ws.Range("J122").Value = 4191.5
$ws.Range("K122").Value = 12600
$ws.Range("L122").Value = 12574.5
$ws.Range("M122").Value = -10150
$ws.Range("N122").Value = -17474.5

$ws.Range("H126").Value = 19303.273
$ws.Range("I126").Value = 4562.4
$ws.Range("J126").Value = 31587.334
$ws.Range("K126").Value = 13687.2
$ws.Range("L126").Value = 94762.00199999999
$ws.Range("M126").Value = -11217.2
$ws.Range("N126").Value = -99702.00199999999

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws.Range("H132").Value = 2728.4707
$ws.Range("I132").Value = 2298.8333
$ws.Range("J132").Value = 3759.6
$ws.Range("K132").Value = 6896.499899999999
$ws.Range("L132").Value = 11278.8
$ws.Range("M132").Value = -4366.499899999999
$ws.Range("N132").Value = -16338.8

$ws.Range("H134").Value = 2636.2856
$ws.Range("I134").Value = 2491.3333
$ws.Range("J134").Value = 2897.2
$ws.Range("K134").Value = 7473.999899999999
$ws.Range("L134").Value = 8691.599999999999
$ws.Range("M134").Value = -4938.999899999999
$ws.Range("N134").Value = -13761.6

$ws.Range("H136").Value = 6091.2104
$ws.Range("I136").Value = 6229.4165
$ws.Range("J136").Value = 5854.2856
$ws.Range("K136").Value = 18688.2495
$ws.Range("L136").Value = 17562.8568
$ws.Range("M136").Value = -16138.2495
$ws.Range("N136").Value = -22662.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 217
$ws.Range("I21").Value = 217
$ws.Range("K21").Value = 651
$ws.Range("M21").Value = -478

$ws.Range("H24").Value = 199.5
$ws.Range("I24").Value = 199.5
$ws.Range("K24").Value = 598.5
$ws.Range("M24").Value = -368.5

$ws.Range("H29").Value = 433.66666
$ws.Range("I29").Value = 11
$ws.Range("J29").Value = 645
$ws.Range("K29").Value = 33
$ws.Range("L29").Value = 1935
$ws.Range("M29").Value = 244
$ws.Range("N29").Value = -2489

$ws.Range("H35").Value = 399
$ws.Range("I35").Value = 399
$ws.Range("K35").Value = 1197
$ws.Range("M35").Value = -909

$ws.Range("H68").Value = 2056.4
$ws.Range("J68").Value = 2157.75
$ws.Range("L68").Value = 6473.25
$ws.Range("N68").Value = -8095.25

$ws.Range("H71").Value = 2056.4
$ws.Range("J71").Value = 2157.75
$ws.Range("L71").Value = 19419.75
$ws.Range("N71").Value = -27531.75

$ws.Range("H107").Value = 1971.5555
$ws.Range("J107").Value = 1992
$ws.Range("L107").Value = 5976
$ws.Range("N107").Value = -9816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 16278.25
$ws.Range("I80").Value = 5066.6
$ws.Range("J80").Value = 24286.572
$ws.Range("K80").Value = 5066.6
$ws.Range("L80").Value = 24286.572
$ws.Range("M80").Value = -4068.6
$ws.Range("N80").Value = -26282.572

$ws.Range("H83").Value = 16278.25
$ws.Range("I83").Value = 5066.6
$ws.Range("J83").Value = 24286.572
$ws.Range("K83").Value = 25333
$ws.Range("L83").Value = 121432.86
$ws.Range("M83").Value = -20341
$ws.Range("N83").Value = -131416.86

$ws.Range("H102").Value = 4168.357
$ws.Range("I102").Value = 4308.5557
$ws.Range("J102").Value = 3916
$ws.Range("K102").Value = 4308.5557
$ws.Range("L102").Value = 3916
$ws.Range("M102").Value = -2686.5557
$ws.Range("N102").Value = -7160

$ws.Range("H126").Value = 4025.7144
$ws.Range("I126").Value = 2239
$ws.Range("K126").Value = 6717
$ws.Range("M126").Value = -4247

$ws.Range("H132").Value = 4178.5454
$ws.Range("I132").Value = 3551.5557
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 10654.6671
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -8124.667099999999
$ws.Range("N132").Value = -26060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4256
$ws.Range("I7").Value = 2756.3333
$ws.Range("J7").Value = 4898.7144
$ws.Range("K7").Value = 2756.3333
$ws.Range("L7").Value = 4898.7144
$ws.Range("M7").Value = -2644.3333
$ws.Range("N7").Value = -5122.7144

$ws.Range("H122").Value = 5678.75
$ws.Range("I122").Value = 4886
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 14658
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -12208
$ws.Range("N122").Value = -25900

$ws.Range("H126").Value = 4256
$ws.Range("I126").Value = 2756.3333
$ws.Range("J126").Value = 4898.7144
$ws.Range("K126").Value = 8268.999899999999
$ws.Range("L126").Value = 14696.1432
$ws.Range("M126").Value = -5798.999899999999
$ws.Range("N126").Value = -19636.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()

$ws.Range("H122").Value = 2730.5173
$ws.Range("I122").Value = 2481.7727
$ws.Range("J122").Value = 3512.2856
$ws.Range("K122").Value = 7445.3181
$ws.Range("L122").Value = 10536.8568
$ws.Range("M122").Value = -4995.3181
$ws.Range("N122").Value = -15436.8568
